{"js": "// Widen the \"Name\" (Nachname) column by 2mm and narrow the \"Vorname\"\n// column by 2mm (130 dxa \u2248 2mm) in the absence table.\n//\n// The table has two adjacent 1418-dxa (70.9pt) columns:\n//   column index 2 (0-based) -> \"Name\"    -> 1548 dxa (77.4pt)\n//   column index 3 (0-based) -> \"Vorname\" -> 1288 dxa (64.4pt)\n//\n// Setting `columnWidth` on a cell resizes the whole column (every cell in\n// that column, plus the <w:gridCol> entry), mirroring what Word does when\n// you drag a column border.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Use the first row (index 0) to address each column once; columnWidth\n// applies to the whole column, not just the single cell.\nconst nameCell = table.getCell(0, 2);\nnameCell.columnWidth = 77.4; // 1548 dxa\n\nconst vornameCell = table.getCell(0, 3);\nvornameCell.columnWidth = 64.4; // 1288 dxa\n\nawait context.sync();\n", "ps1": "# Widen the \"Name\" (Nachname) column by 2mm and narrow the \"Vorname\"\n# column by 2mm (130 dxa ~= 2mm) in the absence table.\n#\n# The table has two adjacent 1418-dxa (70.9pt) columns:\n#   Columns.Item(3) -> \"Name\"    -> 1548 dxa (77.4pt)\n#   Columns.Item(4) -> \"Vorname\" -> 1288 dxa (64.4pt)\n#\n# Word stores widths in twips (dxa) in the OOXML but exposes them to COM\n# in points (1 pt = 20 dxa), so 1548/20 = 77.4 and 1288/20 = 64.4.\n# Setting Column.Width resizes every cell in that column plus the\n# <w:gridCol> entry, same as dragging the column border in the UI.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$nameCol = $t.Columns.Item(3)\n$nameCol.Width = 77.4\n\n$vornameCol = $t.Columns.Item(4)\n$vornameCol.Width = 64.4\n"}
